$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.979.26'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '1.654.10'
$ws.Range("E3").Value = '  +2.21%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'" + '309.74'
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = "'" + '0.3899'
$ws.Range("E7").Value = '  -0.84%  '
$ws.Range("D8").Value = "'" + '0.3835'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = "'" + '51.38'
$ws.Range("E9").Value = '  +4.19%  '
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("D11").Value = "'" + '1.000'
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = "'" + '0.08442'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = "'" + '23.94'
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("D15").Value = "'" + '7.903'
$ws.Range("E15").Value = '  +4.49%  '
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").Value = '1.654.41'
$ws.Range("E17").Value = '  +2.53%  '
$ws.Range("D18").Value = "'" + '94.60'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").Value = "'" + '0.06993'
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").Value = "'" + '19.75'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").Value = "'" + '6.932'
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("D22").Value = "'" + '1.0000'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  +2.04%  '
$ws.Range("D24").Value = '23.981.88'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D25").Value = "'" + '2.464'
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'" + '2.990'
$ws.Range("E26").Value = '  +5.83%  '
$ws.Range("D27").Value = "'" + '22.11'
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("D29").Value = "'" + '5.442'
$ws.Range("E29").Value = '  +2.82%  '
$ws.Range("D30").Value = "'" + '139.01'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").Value = "'" + '7.826'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").Value = "'" + '2.489'
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = '1.835.60'
$ws.Range("E33").Value = '  +2.62%  '
$ws.Range("E34").Value = '  +6.29%  '
$ws.Range("D35").Value = "'" + '0.08103'
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").Value = "'" + '0.02958'
$ws.Range("E36").Value = '  +2.71%  '
$ws.Range("D37").Value = "'" + '6.753'
$ws.Range("E37").Value = '  +2.13%  '
$ws.Range("D38").Value = "'" + '10.92'
$ws.Range("E38").Value = '  +6.04%  '
$ws.Range("D39").Value = "'" + '0.2688'
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = "'" + '0.09145'
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("D41").Value = "'" + '0.7560'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").Value = "'" + '16.44'
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D45").Value = "'" + '0.6948'
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("D46").Value = "'" + '2.458'
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").Value = "'" + '4.092'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").Value = "'" + '0.9998'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = "'" + '0.08290'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").Value = "'" + '1.209'
$ws.Range("E51").Value = '  +0.70%  '
